$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Sertad1"
$ws.Range("C2").Value = "Ar"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 31.82219666666667
$ws.Range("H2").Value = 95.46659
$ws.Range("I2").Value = 0.5609422836697905
$ws.Range("J2").Value = 0.5609422836697905
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 1.119044666666667
$ws.Range("N2").Value = 3.357134
$ws.Range("O2").Value = 0.05243434414602349
$ws.Range("P2").Value = 0.05243434414602349
$ws.Range("Q2").Value = 35.61045946145111
$ws.Range("R2").Value = 320.49413515306
$ws.Range("S2").Value = 0.02941264074799813
$ws.Range("T2").Value = 0.02941264074799813

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Sertad1"
$ws.Range("C3").Value = "Ar"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 31.82219666666667
$ws.Range("H3").Value = 95.46659
$ws.Range("I3").Value = 0.5609422836697905
$ws.Range("J3").Value = 0.5609422836697905
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 15.93829766666667
$ws.Range("N3").Value = 47.814893
$ws.Range("O3").Value = 0.7468103909070326
$ws.Range("P3").Value = 0.7468103909070325
$ws.Range("Q3").Value = 507.1916428805412
$ws.Range("R3").Value = 4564.72478592487
$ws.Range("S3").Value = 0.4189175261437198
$ws.Range("T3").Value = 0.4189175261437197

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Sertad1"
$ws.Range("C4").Value = "Ar"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 31.82219666666667
$ws.Range("H4").Value = 95.46659
$ws.Range("I4").Value = 0.5609422836697905
$ws.Range("J4").Value = 0.5609422836697905
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 4.284484
$ws.Range("N4").Value = 12.853452
$ws.Range("O4").Value = 0.200755264946944
$ws.Range("P4").Value = 0.200755264946944
$ws.Range("Q4").Value = 136.3416924631867
$ws.Range("R4").Value = 1227.07523216868
$ws.Range("S4").Value = 0.1126121167780726
$ws.Range("T4").Value = 0.1126121167780726

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Sertad1"
$ws.Range("C5").Value = "Ar"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 14.52367433333333
$ws.Range("H5").Value = 43.571023
$ws.Range("I5").Value = 0.256014477352223
$ws.Range("J5").Value = 0.256014477352223
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 1.119044666666667
$ws.Range("N5").Value = 3.357134
$ws.Range("O5").Value = 0.05243434414602349
$ws.Range("P5").Value = 0.05243434414602349
$ws.Range("Q5").Value = 16.25264030312022
$ws.Range("R5").Value = 146.273762728082
$ws.Range("S5").Value = 0.0134239512118508
$ws.Range("T5").Value = 0.0134239512118508

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Sertad1"
$ws.Range("C6").Value = "Ar"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 14.52367433333333
$ws.Range("H6").Value = 43.571023
$ws.Range("I6").Value = 0.256014477352223
$ws.Range("J6").Value = 0.256014477352223
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 15.93829766666667
$ws.Range("N6").Value = 47.814893
$ws.Range("O6").Value = 0.7468103909070326
$ws.Range("P6").Value = 0.7468103909070325
$ws.Range("Q6").Value = 231.4826447383932
$ws.Range("R6").Value = 2083.343802645539
$ws.Range("S6").Value = 0.1911942719092734
$ws.Range("T6").Value = 0.1911942719092733

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Sertad1"
$ws.Range("C7").Value = "Ar"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 14.52367433333333
$ws.Range("H7").Value = 43.571023
$ws.Range("I7").Value = 0.256014477352223
$ws.Range("J7").Value = 0.256014477352223
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 4.284484
$ws.Range("N7").Value = 12.853452
$ws.Range("O7").Value = 0.200755264946944
$ws.Range("P7").Value = 0.200755264946944
$ws.Range("Q7").Value = 62.22645030237733
$ws.Range("R7").Value = 560.038052721396
$ws.Range("S7").Value = 0.05139625423109893
$ws.Range("T7").Value = 0.05139625423109893

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Sertad1"
$ws.Range("C8").Value = "Ar"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 10.38402366666667
$ws.Range("H8").Value = 31.152071
$ws.Range("I8").Value = 0.1830432389779865
$ws.Range("J8").Value = 0.1830432389779865
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 1.119044666666667
$ws.Range("N8").Value = 3.357134
$ws.Range("O8").Value = 0.05243434414602349
$ws.Range("P8").Value = 0.05243434414602349
$ws.Range("Q8").Value = 11.62018630272378
$ws.Range("R8").Value = 104.581676724514
$ws.Range("S8").Value = 0.009597752186174565
$ws.Range("T8").Value = 0.009597752186174563

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Sertad1"
$ws.Range("C9").Value = "Ar"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 10.38402366666667
$ws.Range("H9").Value = 31.152071
$ws.Range("I9").Value = 0.1830432389779865
$ws.Range("J9").Value = 0.1830432389779865
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 15.93829766666667
$ws.Range("N9").Value = 47.814893
$ws.Range("O9").Value = 0.7468103909070326
$ws.Range("P9").Value = 0.7468103909070325
$ws.Range("Q9").Value = 165.5036601770448
$ws.Range("R9").Value = 1489.532941593403
$ws.Range("S9").Value = 0.1366985928540395
$ws.Range("T9").Value = 0.1366985928540394

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Sertad1"
$ws.Range("C10").Value = "Ar"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 10.38402366666667
$ws.Range("H10").Value = 31.152071
$ws.Range("I10").Value = 0.1830432389779865
$ws.Range("J10").Value = 0.1830432389779865
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 4.284484
$ws.Range("N10").Value = 12.853452
$ws.Range("O10").Value = 0.200755264946944
$ws.Range("P10").Value = 0.200755264946944
$ws.Range("Q10").Value = 44.49018325545466
$ws.Range("R10").Value = 400.411649299092
$ws.Range("S10").Value = 0.03674689393777246
$ws.Range("T10").Value = 0.03674689393777246

